# Removes leftover decorative shapes that are no longer part of the deck.
#
# Slide 10: drop the full-width white background bar ("object 32") that sat
# behind the slide-number placeholder.
#
# Slide 15: drop the 20 small colored rectangle shapes ("object 57" through
# "object 76") that made up an obsolete highlight grid, leaving the rest of
# the slide (including the slide-number placeholder) untouched.

$p = $ppt.ActivePresentation

$slide10 = $p.Slides.Item(10)
$slide10.Shapes.Item("object 32").Delete()

$slide15 = $p.Slides.Item(15)
for ($i = 57; $i -le 76; $i++) {
    $slide15.Shapes.Item("object $i").Delete()
}
